# Dashboard.pptx edit script
# Reproduces the author's edits:
#  - text/content updates on several slides
#  - delete the old "Summary:" slide and the old blank trailing slide
#  - insert a new "Database" slide after "Data Processing"
#  - insert two duplicate "Results (Christopher)" slides
#  - rename/refill the former "Results" and "Conclusions" slides
#  - refresh the cached footer date (5/22/22 -> 5/23/22) on every layout + master

$p = $ppt.ActivePresentation

function Get-SlideById($id) {
    for ($i = 1; $i -le $p.Slides.Count; $i++) {
        if ($p.Slides.Item($i).SlideID -eq $id) { return $p.Slides.Item($i) }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1. Slide 257 "Introduction:" / "Background:"
# ---------------------------------------------------------------------------
$s257 = Get-SlideById 257
$s257.Shapes.Item(1).TextFrame.TextRange.Text = "Introduction: (Silvania)"
$s257.Shapes.Item(2).TextFrame.TextRange.Text = "Background: "

# ---------------------------------------------------------------------------
# 2. Slide 258 "Data Sources:"
# ---------------------------------------------------------------------------
$s258 = Get-SlideById 258
$s258.Shapes.Item(1).TextFrame.TextRange.Text = "Data Sources: (Kevin) "
$s258.Shapes.Item(2).TextFrame.TextRange.Text = (
    "New York State (countyhealthranking.org)`r" +
    "FIPS Populations Source (census.gov)`r" +
    "`r" +
    "Timeframes"
)

# ---------------------------------------------------------------------------
# 3. Slide 259 "Data Processing:"
# ---------------------------------------------------------------------------
$s259 = Get-SlideById 259
$s259.Shapes.Item(1).TextFrame.TextRange.Text = "Data Processing: (Kevin 2 slides)"
$s259.Shapes.Item(2).TextFrame.TextRange.Text = (
    "EDA etc…`r" +
    "R-Studio used for EDA`r" +
    "Pandas`r" +
    "Jupyter"
)

# ---------------------------------------------------------------------------
# 4. Delete the old "Summary:" slide (263) and the old blank slide (262)
# ---------------------------------------------------------------------------
(Get-SlideById 263).Delete()
(Get-SlideById 262).Delete()

# ---------------------------------------------------------------------------
# 5. Insert new "Database" slide right after slide 259
# ---------------------------------------------------------------------------
$posDb = 0
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).SlideID -eq 259) { $posDb = $i + 1; break }
}
$sDb = $p.Slides.Add($posDb, 2)
$sDb.Shapes.Item(1).TextFrame.TextRange.Text = "Database (Krystal….2 slides)"
$sDb.Shapes.Item(2).TextFrame.TextRange.Text = (
    "PgAdmin & PostGres SQL`r" +
    "Local Database `r" +
    "Connection to DB via Jupyter Notebook for data processing"
)

# ---------------------------------------------------------------------------
# 6. Slide 260 "Machine Learning Models Applied:"
# ---------------------------------------------------------------------------
$s260 = Get-SlideById 260
$s260.Shapes.Item(1).TextFrame.TextRange.Text = "Machine Learning Models Applied: (3 Slides…. Silvania)"
$s260.Shapes.Item(2).TextFrame.TextRange.Text = (
    "Supervised`r" +
    "Unsupervised (PCA & HCA, K-Means)`r" +
    ""
)

# ---------------------------------------------------------------------------
# 7. Slide 261 "Results" -> "Results (Christopher)"
# ---------------------------------------------------------------------------
$resultsBody = (
    "Visualizations `r" +
    "Embed Tableau link in powerpoint`r" +
    "2 or slides to summarise key findings (full details in tableau)`r" +
    ""
)

$s261 = Get-SlideById 261
$s261.Shapes.Item(1).TextFrame.TextRange.Text = "Results (Christopher)"
$s261.Shapes.Item(2).TextFrame.TextRange.Text = $resultsBody

# ---------------------------------------------------------------------------
# 8. Insert two duplicate "Results (Christopher)" slides after slide 261
# ---------------------------------------------------------------------------
$posRes = 0
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).SlideID -eq 261) { $posRes = $i + 1; break }
}
$sR1 = $p.Slides.Add($posRes, 2)
$sR1.Shapes.Item(1).TextFrame.TextRange.Text = "Results (Christopher)"
$sR1.Shapes.Item(2).TextFrame.TextRange.Text = $resultsBody

$sR2 = $p.Slides.Add($posRes + 1, 2)
$sR2.Shapes.Item(1).TextFrame.TextRange.Text = "Results (Christopher)"
$sR2.Shapes.Item(2).TextFrame.TextRange.Text = $resultsBody

# ---------------------------------------------------------------------------
# 9. Slide 264 "Conclusions and Recommendations" -> move to end & refill
# ---------------------------------------------------------------------------
$s264 = Get-SlideById 264
$s264.Shapes.Item(1).TextFrame.TextRange.Text = "Conclusions and Recommendations (Everybody)"
$s264.Shapes.Item(2).TextFrame.TextRange.Text = "Suggestions for refinements or further areas of study"
$s264.MoveTo($p.Slides.Count)

# ---------------------------------------------------------------------------
# 10. Refresh the cached footer date on every layout + the slide master
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.PlaceholderFormat.Type -eq 16) {
        $sh.TextFrame.TextRange.Text = "5/23/22"
    }
}
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $cl = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        if ($sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = "5/23/22"
        }
    }
}

Write-Output "Final slide count: $($p.Slides.Count)"
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    Write-Output "$i : id=$($sl.SlideID) title=$($sl.Shapes.Item(1).TextFrame.TextRange.Text)"
}
